$d = $word.ActiveDocument

# 1. Update the delivery date in the header
$d.Content.Find.Execute("2011-11-03", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2011-11-25", 2)

# 2. Merge "se  cuenta con " + "2" runs into a single run "se  cuenta con 2"
$d.Content.Find.Execute("se  cuenta con " + [char]2, $false, $false, $false, $false, $false,
                         $true, 1, $false, "se  cuenta con 2", 2)

# 3. Merge "mayor a " + "2" runs into a single run "mayor a 2"
$d.Content.Find.Execute("mayor a " + [char]2, $false, $false, $false, $false, $false,
                         $true, 1, $false, "mayor a 2", 2)
